$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff items produced by this report-generation run:
#   - 3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png  (existing dependency, re-handed-off)
#   - 8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md    (new source file)
#   - ada3ae2c-047b-4b68-b4e9-584b943a06a1.png   (new dependency of the .md file)
# ---------------------------------------------------------------------------

# =================================== Overview ===============================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("B2:C2").Value2 = "Ready for handoff"
$ws1.Range("D2").Value2 = "2016-49-17 14:49:17"

$arr = New-Object 'object[,]' 2,3
$arr[0,0] = "Ready for handoff"; $arr[0,1] = "Ready for handoff"; $arr[0,2] = "2016-49-17 14:49:17"
$arr[1,0] = "Ready for handoff"; $arr[1,1] = "Ready for handoff"; $arr[1,2] = "2016-49-17 14:49:17"
$ws1.Range("B3:D4").Value2 = $arr

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png", [Type]::Missing, [Type]::Missing, "3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md", [Type]::Missing, [Type]::Missing, "8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/ada3ae2c-047b-4b68-b4e9-584b943a06a1.png", [Type]::Missing, [Type]::Missing, "ada3ae2c-047b-4b68-b4e9-584b943a06a1.png") | Out-Null

# =================================== zh-cn ===================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("C2").Value2 = "Ready for handoff"
$ws2.Range("E2").Value2 = "2016-03-17 14:49:12"
$ws2.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I2").Value2 = "IsDependency"
$ws2.Range("J2").Value2 = "e2e\8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md"

$ws2.Range("C3:C3").Value2 = "Ready for handoff"
$ws2.Range("E3").Value2 = "2016-03-17 14:49:12"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I3").Value2 = "Include"

$ws2.Range("C4").Value2 = "Ready for handoff"
$ws2.Range("E4").Value2 = "2016-03-17 14:49:12"
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("I4").Value2 = "IsDependency"
$ws2.Range("J4").Value2 = "e2e\8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png", [Type]::Missing, [Type]::Missing, "3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png", [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25feb1ce67db86d0cafcf918e890ca710617bcd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d6c43709b246e490165bcef7f8a471fd046ce842.png", [Type]::Missing, [Type]::Missing, "d6c43709b246e490165bcef7f8a471fd046ce842.png") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md", [Type]::Missing, [Type]::Missing, "8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25feb1ce67db86d0cafcf918e890ca710617bcd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.d948b8b82ea1d7550c421631c16f0ff0a1be5c16.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8ade00eb-d889-4cfa-a80a-2081ab0cec3b.d948b8b82ea1d7550c421631c16f0ff0a1be5c16.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/ada3ae2c-047b-4b68-b4e9-584b943a06a1.png", [Type]::Missing, [Type]::Missing, "ada3ae2c-047b-4b68-b4e9-584b943a06a1.png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/ada3ae2c-047b-4b68-b4e9-584b943a06a1.png", [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25feb1ce67db86d0cafcf918e890ca710617bcd4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f017df5f14adca339854985eb4574d9f1ce11c5f.png", [Type]::Missing, [Type]::Missing, "f017df5f14adca339854985eb4574d9f1ce11c5f.png") | Out-Null

# =================================== de-de ===================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("C2").Value2 = "Ready for handoff"
$ws3.Range("E2").Value2 = "2016-03-17 14:49:17"
$ws3.Range("H2").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I2").Value2 = "IsDependency"
$ws3.Range("J2").Value2 = "e2e\8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md"

$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("E3").Value2 = "2016-03-17 14:49:17"
$ws3.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I3").Value2 = "Include"

$ws3.Range("C4").Value2 = "Ready for handoff"
$ws3.Range("E4").Value2 = "2016-03-17 14:49:17"
$ws3.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("I4").Value2 = "IsDependency"
$ws3.Range("J4").Value2 = "e2e\8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png", [Type]::Missing, [Type]::Missing, "3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/3ab69aa5-cd61-43eb-9ce7-54a30bea6118.png", [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cee03e04526835a03268df846267b6bd143a84a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d6c43709b246e490165bcef7f8a471fd046ce842.png", [Type]::Missing, [Type]::Missing, "d6c43709b246e490165bcef7f8a471fd046ce842.png") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md", [Type]::Missing, [Type]::Missing, "8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md", [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cee03e04526835a03268df846267b6bd143a84a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8ade00eb-d889-4cfa-a80a-2081ab0cec3b.d948b8b82ea1d7550c421631c16f0ff0a1be5c16.de-de.xlf", [Type]::Missing, [Type]::Missing, "8ade00eb-d889-4cfa-a80a-2081ab0cec3b.d948b8b82ea1d7550c421631c16f0ff0a1be5c16.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/ada3ae2c-047b-4b68-b4e9-584b943a06a1.png", [Type]::Missing, [Type]::Missing, "ada3ae2c-047b-4b68-b4e9-584b943a06a1.png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/04842f7cc1659a3ee549d305307f35eeb32223f2/e2e/ada3ae2c-047b-4b68-b4e9-584b943a06a1.png", [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cee03e04526835a03268df846267b6bd143a84a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f017df5f14adca339854985eb4574d9f1ce11c5f.png", [Type]::Missing, [Type]::Missing, "f017df5f14adca339854985eb4574d9f1ce11c5f.png") | Out-Null

Write-Host "Handback report regenerated: added rows for 8ade00eb-d889-4cfa-a80a-2081ab0cec3b.md and its dependencies."
